$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("C2").Value = 0.1513560054718786
$ws.Range("D2").Value = -6.294687658057581
$ws.Range("E2").Value = -6.294687658057581
$ws.Range("F2").Value = -0.8269167698135087
$ws.Range("I2").Value = 10.2512112733724
$ws.Range("J2").Value = -0.1524511281828686
$ws.Range("L2").Value = -0.2622156542051112
$ws.Range("M2").Value = 8.566247266248954

# Row 32 updates
$ws.Range("C32").Value = -0.5024491423859256
$ws.Range("F32").Value = -1.348648551689004
$ws.Range("G32").Value = 5.534708212411552
$ws.Range("J32").Value = -0.8213924077415695
$ws.Range("L32").Value = -1
